# Move the 3 "new" listings from the "New" sheet onto the end of the
# "Previously added" sheet (as rows 219-221), then leave "New" with only
# its header row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

$lastRow = 218
$firstNewRow = $lastRow + 1     # 219
$moveCount = 3                   # rows 2..4 on "New"
$lastNewRow = $firstNewRow + $moveCount - 1   # 221

# 1) Stamp the 3 destination rows with the same look & feel (styles/number
#    formats) as the current last row of data, by copying its formatting.
$ws1.Range("A" + $lastRow + ":F" + $lastRow).Copy()
$ws1.Range("A" + $firstNewRow + ":F" + $lastNewRow).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 2) Bring over the actual data (values only, so underlying types such as
#    numeric-looking text stay text) row by row from "New" -> "Previously added".
for ($i = 0; $i -lt $moveCount; $i++) {
    $srcRow = 2 + $i
    $dstRow = $firstNewRow + $i
    $ws2.Range("A" + $srcRow + ":F" + $srcRow).Copy()
    $ws1.Range("A" + $dstRow + ":F" + $dstRow).PasteSpecial(-4163)      # xlPasteValues
}
$excel.CutCopyMode = 0

# 3) Recreate the column-A hyperlinks (the link target is simply the URL
#    text that's already in the cell) on the destination rows.
for ($i = 0; $i -lt $moveCount; $i++) {
    $dstRow = $firstNewRow + $i
    $cell = $ws1.Cells.Item($dstRow, 1)
    $url = [string]$cell.Value2
    $ws1.Hyperlinks.Add($cell, $url)
}

# Adding hyperlinks re-applies Excel's built-in "Hyperlink" cell style, so
# restore the original column-A formatting that was set up in step 1.
$ws1.Range("A" + $lastRow).Copy()
$ws1.Range("A" + $firstNewRow + ":A" + $lastNewRow).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ... and drop the now-unused "Hyperlink" named cell style that Excel
# auto-registered in the workbook when the links were added.
try {
    $wb.Styles.Item("Hyperlink").Delete()
} catch {
}

# 4) Clean out the rows (and their hyperlinks) from "New", leaving just
#    the header row behind.
$ws2.Hyperlinks.Delete()
$ws2.Range("A2:F" + ($moveCount + 1)).Delete(-4162)   # xlShiftUp
